$d = $word.ActiveDocument

function Split-Run($startPos, $endPos) {
    # Touch Bold (set then restore) on the inserted range so the serializer is
    # forced to break the run at both of its boundaries instead of silently
    # re-merging text that has identical formatting on every side of the split.
    $probe = $d.Range($startPos, $endPos)
    if ($probe.End -gt $probe.Start) {
        $was = $probe.Font.Bold
        $probe.Font.Bold = 1
        $probe.Font.Bold = $was
    }
}

# ---------------------------------------------------------------------------
# 1) "...open the c:\wxWidgets\build\msw\wx_vc9.sln file"
#    -> insert "32" right after "wxWidgets" (its own run), then move the
#       _GoBack bookmark to sit immediately after the "32", right before
#       "\build...". Bookmarks.Add with an existing name relocates it, so
#       this also clears it out of the trailing empty paragraph it used to
#       occupy.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Open Visual Studio 2012, and open the c:\wxWidgets\build\msw\wx_vc9.sln file", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $phraseStart = $rng.Start
    $insertPos = $phraseStart + "Open Visual Studio 2012, and open the c:\wxWidgets".Length

    $ip = $d.Range($insertPos, $insertPos)
    $ip.InsertAfter("32")

    # Force "32" onto its own run, separate from the preceding text.
    $insertEnd = $insertPos + 2
    Split-Run $insertPos $insertEnd

    # Bookmark sits right after "32", before "\build".
    $bmPos = $insertPos + 2
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# 2) "Set the WXMSW3 environment variable to c:\wxWidgets" -> append "32" as
#    its own trailing run.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Set the WXMSW3 environment variable to c:\wxWidgets", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $endPos = $rng2.End
    $ip2 = $d.Range($endPos, $endPos)
    $ip2.InsertAfter("32")
    $endPos2 = $endPos + 2
    Split-Run $endPos $endPos2
}

Write-Output "ok"
